# Minor design changes to the ToolTip Dialog
# - Adds 3 new translation rows (Percentage/Prozentwert, Region, Info) to each
#   of the three sheets (Exiobase, Deutsch, English).
# - Updates selection / active-sheet view state to match where the author
#   ended up after editing (English sheet active, scrolled further down).

$wb = $excel.ActiveWorkbook

$wsExio = $wb.Worksheets.Item("Exiobase")
$wsDe   = $wb.Worksheets.Item("Deutsch")
$wsEn   = $wb.Worksheets.Item("English")

# ---------------------------------------------------------------------
# New shared strings must be minted in this exact order so they land on
# the same indices as the target file: 117 Percentage, 118 Prozentwert,
# 119 Region, 120 Info.
# ---------------------------------------------------------------------
$wsExio.Range("A59").Value = "Percentage"   # mints shared string 117
$wsDe.Range("B59").Value   = "Prozentwert"  # mints shared string 118
$wsExio.Range("A60").Value = "Region"       # mints shared string 119
$wsExio.Range("A61").Value = "Info"         # mints shared string 120

# ---------------------------------------------------------------------
# Exiobase sheet (English key == English key, rows 59-61)
# ---------------------------------------------------------------------
$wsExio.Range("B59").Value = "Percentage"

$wsExio.Range("B60").Value = "Region"

$wsExio.Range("B61").Value = "Info"

$wsExio.Range("A59:B61").RowHeight = 15
$wsExio.Range("A59:B61").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# Deutsch sheet (English key -> German translation, rows 59-61)
# ---------------------------------------------------------------------
$wsDe.Range("A59").Value = "Percentage"

$wsDe.Range("A60").Value = "Region"
$wsDe.Range("B60").Value = "Region"

$wsDe.Range("A61").Value = "Info"
$wsDe.Range("B61").Value = "Info"

$wsDe.Range("A59:B61").RowHeight = 15
$wsDe.Range("A59:B61").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# English sheet (English key == English key, rows 63-65)
# ---------------------------------------------------------------------
$wsEn.Range("A63").Value = "Percentage"
$wsEn.Range("B63").Value = "Percentage"

$wsEn.Range("A64").Value = "Region"
$wsEn.Range("B64").Value = "Region"

$wsEn.Range("A65").Value = "Info"
$wsEn.Range("B65").Value = "Info"

$wsEn.Range("A63:B65").RowHeight = 15
$wsEn.Range("A63:B65").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# View / selection state: select in order Exiobase -> Deutsch -> English so
# English ends up the active tab (matches the saved activeTab/tabSelected).
# ---------------------------------------------------------------------
$wsExio.Range("A61").Select()
$wsDe.Range("B61").Select()
$wsEn.Range("B65").Select()

Write-Output "done"
